$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column G width (OOXML width="17")
$ws.Columns.Item(7).ColumnWidth = 16.1667

# Insert a new row 9 for Yellow LED (shifts old row9 Green LED -> row10, old row10 Slide switch -> row11)
$ws.Rows.Item(9).Insert()

# New row 12: Reset switch (new line item)
$ws.Range("B12").Value = 1
$ws.Range("G12").Value = "Reset switch"
$ws.Range("E12").Value = "EG5384CT-ND"
$ws.Range("D12").Value = "TL3365AF180QG"
$ws.Range("C12").Value = "SWITCH TACTILE SPST-NO 0.05A 12V"

# Row 11: Slide switch (shifted down), add note
$ws.Range("G11").Value = "5V power to target switch"

# Row 7: 3.3V power regulator
$ws.Range("G7").Value = "3.3V power regulator"

# Row 6: CSI debug/programming male header
$ws.Range("G6").Value = "CSI debug/programming male header"

# Row 4: add MCU designator note in G
$ws.Range("G4").Value = "MCU"

# Row 13: .1uF cap note only
$ws.Range("G13").Value = ".1uF cap"

# Row 14: 4.7K resistor, quantity 2
$ws.Range("B14").Value = 2
$ws.Range("G14").Value = "4.7K resistor"

# Row 9: Yellow LED
$ws.Range("C9").Value = "Yellow LED"

# Row 8: Red LED -> ST-Link LED note
$ws.Range("G8").Value = "ST-Link LED"
$ws.Range("G9").Value = "ST-Link LED"

# Row 10: Green LED (shifted down), add note
$ws.Range("G10").Value = "3.3V target power LED"

# Update selection to match target
$ws.Range("G11").Select()
